$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) / Volume(1h) (E) columns store plain text values in the
# source data (e.g. "257.37", "0.99%"), not numbers. Force each touched
# cell to Text format before assigning so Excel keeps the literal string
# instead of auto-converting it to a number/percentage.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"

$ws.Range("D2").Value = "257.37"
$ws.Range("E2").Value = "0.99%"
$ws.Range("D3").Value = "27.04"
$ws.Range("E3").Value = "-4.06%"
$ws.Range("D4").Value = "4.764"
$ws.Range("E4").Value = "-10.58%"
$ws.Range("D5").Value = "0.05966"
$ws.Range("E5").Value = "2.06%"
$ws.Range("D6").Value = "6.677"
$ws.Range("E6").Value = "-0.47%"
$ws.Range("D7").Value = "0.8703"
$ws.Range("E7").Value = "0.46%"
$ws.Range("D8").Value = "0.9451"
$ws.Range("E8").Value = "3.65%"
$ws.Range("D9").Value = "0.1405"
$ws.Range("E9").Value = "-1.17%"
$ws.Range("D10").Value = "0.03613"
$ws.Range("E10").Value = "4.83%"
$ws.Range("D11").Value = "0.07169"
$ws.Range("E11").Value = "0.25%"
$ws.Range("D12").Value = "0.03167"
$ws.Range("E12").Value = "-0.39%"
$ws.Range("D13").Value = "0.09239"
$ws.Range("E13").Value = "0.23%"
$ws.Range("D14").Value = "0.001548"
$ws.Range("E14").Value = "0.71%"
$ws.Range("D15").Value = "0.0006111"
$ws.Range("E15").Value = "0.88%"
$ws.Range("D16").Value = "0.005995"
$ws.Range("E16").Value = "3.05%"
$ws.Range("D17").Value = "3.474"
$ws.Range("E17").Value = "-0.69%"
$ws.Range("D18").Value = "3.176"
$ws.Range("E18").Value = "-1.70%"
$ws.Range("E19").Value = "1.76%"
$ws.Range("D20").Value = "0.3100"
$ws.Range("E20").Value = "-2.26%"
$ws.Range("D21").Value = "0.1306"
$ws.Range("E21").Value = "-0.70%"
$ws.Range("D22").Value = "3.814"
$ws.Range("E22").Value = "7.73%"
$ws.Range("D23").Value = "0.04222"
$ws.Range("E23").Value = "1.35%"
$ws.Range("D25").Value = "0.001222"
$ws.Range("E25").Value = "-0.43%"
$ws.Range("D26").Value = "0.004499"
$ws.Range("E26").Value = "-10.76%"
$ws.Range("D27").Value = "0.0001200"
$ws.Range("D28").Value = "0.0001493"
$ws.Range("E28").Value = "-22.92%"
$ws.Range("D40").Value = "0.03819"
$ws.Range("E40").Value = "-0.87%"
$ws.Range("D41").Value = "0.006171"
$ws.Range("E41").Value = "7.80%"
$ws.Range("D42").Value = "0.1101"
$ws.Range("E42").Value = "-0.01%"
$ws.Range("D43").Value = "0.002252"
$ws.Range("E43").Value = "2.45%"
$ws.Range("D44").Value = "0.01057"
$ws.Range("E44").Value = "-3.77%"
$ws.Range("D45").Value = "0.00005500"
$ws.Range("E45").Value = "5.18%"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").Value = "0.10%"
$ws.Range("D47").Value = "0.1091"
$ws.Range("E47").Value = "21.49%"
$ws.Range("E48").Value = "5.63%"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").Value = "0.10%"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").Value = "0.10%"
